# Update the cryptos list: refresh Price (D) and Volume(1h) (E) columns.
# D/E cells hold plain text (e.g. "52.303.12", "  -0.05%  "), not real
# numbers/percentages, so values are written as ="literal" formulas first
# (this avoids Excel's automatic text->number coercion that would occur
# with a plain .Value assignment) and then converted to static text via
# PasteSpecial(xlPasteValues) in one pass, matching the original data type
# and leaving cell styling untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="52.317.84"'
$ws.Range('E2').Formula = '="  -0.01%  "'
$ws.Range('D3').Formula = '="2.839.84"'
$ws.Range('E3').Formula = '="  +1.67%  "'
$ws.Range('D5').Formula = '="361.06"'
$ws.Range('E5').Formula = '="  +3.95%  "'
$ws.Range('D6').Formula = '="112.98"'
$ws.Range('E6').Formula = '="  -2.50%  "'
$ws.Range('E7').Formula = '="  +4.11%  "'
$ws.Range('E8').Formula = '="  +0.01%  "'
$ws.Range('D9').Formula = '="0.604"'
$ws.Range('E9').Formula = '="  +2.23%  "'
$ws.Range('D10').Formula = '="41.23"'
$ws.Range('E10').Formula = '="  -2.97%  "'
$ws.Range('D11').Formula = '="0.0881"'
$ws.Range('E11').Formula = '="  +2.47%  "'
$ws.Range('E12').Formula = '="  +0.92%  "'
$ws.Range('D13').Formula = '="20.13"'
$ws.Range('E13').Formula = '="  +0.64%  "'
$ws.Range('D14').Formula = '="7.84"'
$ws.Range('E14').Formula = '="  -0.39%  "'
$ws.Range('D15').Formula = '="3.286.48"'
$ws.Range('E15').Formula = '="  +1.78%  "'
$ws.Range('D16').Formula = '="2.856.31"'
$ws.Range('E16').Formula = '="  +2.57%  "'
$ws.Range('D17').Formula = '="0.934"'
$ws.Range('E17').Formula = '="  +4.84%  "'
$ws.Range('D18').Formula = '="52.240.98"'
$ws.Range('E18').Formula = '="  +0.02%  "'
$ws.Range('E19').Formula = '="  +3.92%  "'
$ws.Range('D20').Formula = '="3.16"'
$ws.Range('E20').Formula = '="  -0.65%  "'
$ws.Range('D21').Formula = '="13.56"'
$ws.Range('E21').Formula = '="  +1.70%  "'
$ws.Range('E22').Formula = '="  +2.37%  "'
$ws.Range('D23').Formula = '="273.12"'
$ws.Range('E23').Formula = '="  +1.26%  "'
$ws.Range('D24').Formula = '="70.72"'
$ws.Range('E24').Formula = '="  +0.95%  "'
$ws.Range('E25').Formula = '="  +3.19%  "'
$ws.Range('D26').Formula = '="27.27"'
$ws.Range('E26').Formula = '="  +1.73%  "'
$ws.Range('D28').Formula = '="10.39"'
$ws.Range('E28').Formula = '="  +1.36%  "'
$ws.Range('E29').Formula = '="  +0.27%  "'
$ws.Range('E30').Formula = '="  +2.41%  "'
$ws.Range('D31').Formula = '="0.0484"'
$ws.Range('E31').Formula = '="  +6.49%  "'
$ws.Range('D32').Formula = '="35.58"'
$ws.Range('E32').Formula = '="  +3.37%  "'
$ws.Range('D33').Formula = '="52.35"'
$ws.Range('E33').Formula = '="  +4.42%  "'
$ws.Range('E34').Formula = '="  +3.07%  "'
$ws.Range('D35').Formula = '="5.65"'
$ws.Range('E35').Formula = '="  +14.31%  "'
$ws.Range('E36').Formula = '="  +3.07%  "'
$ws.Range('E37').Formula = '="  -0.11%  "'
$ws.Range('D38').Formula = '="3.31"'
$ws.Range('E38').Formula = '="  +2.95%  "'
$ws.Range('E39').Formula = '="  -2.45%  "'
$ws.Range('D40').Formula = '="18.57"'
$ws.Range('E40').Formula = '="  -0.12%  "'
$ws.Range('E41').Formula = '="  +2.35%  "'
$ws.Range('D42').Formula = '="127.23"'
$ws.Range('E42').Formula = '="  +0.37%  "'
$ws.Range('D43').Formula = '="2.55"'
$ws.Range('E43').Formula = '="  -1.93%  "'
$ws.Range('D44').Formula = '="23.34"'
$ws.Range('E44').Formula = '="  +0.33%  "'
$ws.Range('D45').Formula = '="2.28"'
$ws.Range('E45').Formula = '="  -0.42%  "'
$ws.Range('E46').Formula = '="  +2.20%  "'
$ws.Range('D47').Formula = '="2.098.32"'
$ws.Range('E47').Formula = '="  +2.03%  "'
$ws.Range('E48').Formula = '="  -1.38%  "'
$ws.Range('D49').Formula = '="5.93"'
$ws.Range('E49').Formula = '="  +5.58%  "'
$ws.Range('D50').Formula = '="0.979"'
$ws.Range('E50').Formula = '="  +1.53%  "'
$ws.Range('D51').Formula = '="9.26"'
$ws.Range('E51').Formula = '="  +3.27%  "'

$rng = $ws.Range("D2:E51")
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
